# Generate Report for Archive
# Update the localization-status report: a22b1761 moves to "In Translation",
# and the ab566d84 / cd794114 rows swap places (row 4 <-> row 5) on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# a22b1761 row (row 3): status changes to "In Translation"
$ws.Range("B3").Value = "In Translation"
$ws.Range("C3").Value = "In Translation"

# Row 4 becomes cd794114's data, row 5 becomes ab566d84's data (swap)
$ws.Range("A4").Value = "cd794114-f34c-4c9a-b442-8dc4ca874ae7.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "2016-03-23 02:37:47"

$ws.Range("A5").Value = "ab566d84-bc18-4cb6-98e0-d2e8d09b4db5.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "2016-03-23 02:36:09"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$4') { $hl.TextToDisplay = "cd794114-f34c-4c9a-b442-8dc4ca874ae7.md" }
    elseif ($addr -eq '$A$5') { $hl.TextToDisplay = "ab566d84-bc18-4cb6-98e0-d2e8d09b4db5.md" }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# a22b1761 row (row 3): status changes to "In Translation"
$ws.Range("C3").Value = "In Translation"

# Row 4 becomes cd794114's data, row 5 becomes ab566d84's data (swap)
$ws.Range("A4").Value = "cd794114-f34c-4c9a-b442-8dc4ca874ae7.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "cd794114-f34c-4c9a-b442-8dc4ca874ae7.f1061303dae04b5d93af818dd52789641ba024d4.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-23 02:37:43"

$ws.Range("A5").Value = "ab566d84-bc18-4cb6-98e0-d2e8d09b4db5.md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "ab566d84-bc18-4cb6-98e0-d2e8d09b4db5.085fd1f70cf3af61292510ce665b16cedb4c0f3f.zh-cn.xlf"
$ws.Range("E5").Value = "2016-03-23 02:36:05"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$4') { $hl.TextToDisplay = "cd794114-f34c-4c9a-b442-8dc4ca874ae7.md" }
    elseif ($addr -eq '$D$4') { $hl.TextToDisplay = "cd794114-f34c-4c9a-b442-8dc4ca874ae7.f1061303dae04b5d93af818dd52789641ba024d4.zh-cn.xlf" }
    elseif ($addr -eq '$A$5') { $hl.TextToDisplay = "ab566d84-bc18-4cb6-98e0-d2e8d09b4db5.md" }
    elseif ($addr -eq '$D$5') { $hl.TextToDisplay = "ab566d84-bc18-4cb6-98e0-d2e8d09b4db5.085fd1f70cf3af61292510ce665b16cedb4c0f3f.zh-cn.xlf" }
}

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# a22b1761 row (row 3): status changes to "In Translation"
$ws.Range("C3").Value = "In Translation"

# Row 4 becomes cd794114's data, row 5 becomes ab566d84's data (swap)
$ws.Range("A4").Value = "cd794114-f34c-4c9a-b442-8dc4ca874ae7.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "cd794114-f34c-4c9a-b442-8dc4ca874ae7.f1061303dae04b5d93af818dd52789641ba024d4.de-de.xlf"
$ws.Range("E4").Value = "2016-03-23 02:37:47"

$ws.Range("A5").Value = "ab566d84-bc18-4cb6-98e0-d2e8d09b4db5.md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "ab566d84-bc18-4cb6-98e0-d2e8d09b4db5.085fd1f70cf3af61292510ce665b16cedb4c0f3f.de-de.xlf"
$ws.Range("E5").Value = "2016-03-23 02:36:09"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$4') { $hl.TextToDisplay = "cd794114-f34c-4c9a-b442-8dc4ca874ae7.md" }
    elseif ($addr -eq '$D$4') { $hl.TextToDisplay = "cd794114-f34c-4c9a-b442-8dc4ca874ae7.f1061303dae04b5d93af818dd52789641ba024d4.de-de.xlf" }
    elseif ($addr -eq '$A$5') { $hl.TextToDisplay = "ab566d84-bc18-4cb6-98e0-d2e8d09b4db5.md" }
    elseif ($addr -eq '$D$5') { $hl.TextToDisplay = "ab566d84-bc18-4cb6-98e0-d2e8d09b4db5.085fd1f70cf3af61292510ce665b16cedb4c0f3f.de-de.xlf" }
}
